# Rename the worksheet to reflect the new "last updated" date
# (02-12-2025 -> 05-12-2025). Excel automatically propagates the new
# sheet name into any defined names / formulas that reference the
# sheet (e.g. the "Fysioterapisystemer" defined name range), so no
# further action is required for that part of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Opdateret d. 05-12-2025"
